$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark left in the title line.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2. The "»»»  your stuff after this line »»»" paragraph is split across
#    three runs (with gramStart/gramEnd proof-error markers around
#    "»  your"). Re-running Find/Replace over the exact same visible text
#    collapses it back down to a single plain run.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    ">>>  your stuff after this line >>>", $true, $false, $false, $false,
    $false, $true, 1, $false, ">>>  your stuff after this line >>>", 2
) | Out-Null

# ---------------------------------------------------------------------
# 3. Insert a new paragraph right after that line with red text reading
#    "This File is changed."
# ---------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like ">>>  your stuff after this line >>>*") {
        $target = $p
    }
}

if ($target -ne $null) {
    $target.Range.InsertParagraphAfter()
    $newPara = $target.Next()
    $newPara.Range.Text = "This File is changed."
    $newPara.Range.Font.Color = 255
}
